$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly schedule is being finished: time slots 14:00-16:00 are split into
# 12-minute sub-blocks, several day/category assignments change, and the last
# "Familia y descanso" row is duplicated onto a new 22:00-23:00 row.
# Start from a clean slate over the old table and rewrite it completely.
$ws.Range("A1:H16").ClearContents()

# Row 1
$ws.Range("A1").Value = "Horas"
$ws.Range("B1").Value = "Lunes"
$ws.Range("C1").Value = "Martes"
$ws.Range("D1").Value = "Miercoles"
$ws.Range("E1").Value = "Jueves"
$ws.Range("F1").Value = "Viernes"
$ws.Range("G1").Value = "Sabado"
$ws.Range("H1").Value = "Domingo"

# Row 2
$ws.Range("A2").Value = "8:00 - 9:00"

# Row 3
$ws.Range("A3").Value = "9:00 - 10:00"
$ws.Range("B3").Value = "Clases"
$ws.Range("C3").Value = "Clases"
$ws.Range("D3").Value = "Clases"
$ws.Range("E3").Value = "Clases"
$ws.Range("F3").Value = "Clases"
$ws.Range("G3").Value = "Clases"

# Row 4
$ws.Range("A4").Value = "10:00 - 11:00"
$ws.Range("B4").Value = "Clases"
$ws.Range("C4").Value = "Clases"
$ws.Range("D4").Value = "Clases"
$ws.Range("E4").Value = "Clases"
$ws.Range("F4").Value = "Clases"
$ws.Range("G4").Value = "Clases"

# Row 5
$ws.Range("A5").Value = "11:00 - 12:00"
$ws.Range("F5").Value = "Ocio"
$ws.Range("G5").Value = "Ocio"

# Row 6
$ws.Range("A6").Value = "12:00 - 13:00"
$ws.Range("F6").Value = "Ocio"
$ws.Range("G6").Value = "Ocio"

# Row 7
$ws.Range("A7").Value = "13:00 - 14:00"
$ws.Range("F7").Value = "Estudio"

# Row 8
$ws.Range("A8").Value = "14:00 - 14:12"
$ws.Range("F8").Value = "Descanso"

# Row 9
$ws.Range("A9").Value = "14:12 - 15:00"
$ws.Range("F9").Value = "Estudio"

# Row 10
$ws.Range("A10").Value = "15:00 - 15:12"

# Row 11
$ws.Range("A11").Value = "15:12 - 16:00"

# Row 12
$ws.Range("A12").Value = "16:00 - 17:00"

# Row 13
$ws.Range("A13").Value = "17:00 - 18:00"

# Row 14
$ws.Range("A14").Value = "18:00 - 19:00"

# Row 15
$ws.Range("A15").Value = "19:00 - 20:00"

# Row 16
$ws.Range("A16").Value = "20:00 - 21:00"

# Row 17
$ws.Range("A17").Value = "21:00 - 22:00"
$ws.Range("B17").Value = "Familia y descanso"
$ws.Range("C17").Value = "Familia y descanso"
$ws.Range("D17").Value = "Familia y descanso"
$ws.Range("E17").Value = "Familia y descanso"
$ws.Range("F17").Value = "Familia y descanso"
$ws.Range("G17").Value = "Familia y descanso"

# Row 18
$ws.Range("A18").Value = "22:00 - 23:00"
$ws.Range("B18").Value = "Familia y descanso"
$ws.Range("C18").Value = "Familia y descanso"
$ws.Range("D18").Value = "Familia y descanso"
$ws.Range("E18").Value = "Familia y descanso"
$ws.Range("F18").Value = "Familia y descanso"
$ws.Range("G18").Value = "Familia y descanso"
